$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 21): Date, Total Count, Session Timeout Errors, Errors Requiring Analysis
$ws.Range("A21").Value = 45965
$ws.Range("B21").Value = 716
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 694

# Move the active selection to the newly added row, matching the saved UI state
$excel.Goto($ws.Range("A21:D21"))
